# Apply the authored edit to ISD_Project_Data.xlsx:
#  - Remove the floating "TextBox 1" shape (math formula ΔD_q^p) from the
#    Change_Demand sheet and replace it with an equivalent plain-text label
#    "q\p" typed directly into cell A1 of that sheet.
#  - Zero out four cells in the Change_Demand matrix (H3, I5, C8, E9).
#  - Make Change_Demand the active sheet/tab (previously Recapture_Probability
#    was active), and update the remembered selections on the Itenaries and
#    Change_Demand sheets to E12.

$wb = $excel.ActiveWorkbook

$wsItenaries = $wb.Worksheets.Item("Itenaries")
$wsChangeDemand = $wb.Worksheets.Item("Change_Demand")

# Delete the floating text-box shape that rendered the math label on the
# Change_Demand sheet; its content is being replaced by a plain cell value.
foreach ($shp in $wsChangeDemand.Shapes) {
    [void]$shp.Delete()
}

# The shape used to display the label "q\p" (math notation for ΔD with sub
# q and sup p) above the matrix; put the same text directly into A1 now.
$wsChangeDemand.Range("A1").Value = "q\p"

# Update the matrix values that changed.
$wsChangeDemand.Range("H3").Value = 0
$wsChangeDemand.Range("I5").Value = 0
$wsChangeDemand.Range("C8").Value = 0
$wsChangeDemand.Range("E9").Value = 0

# Update remembered selections.
[void]$wsItenaries.Activate()
[void]$wsItenaries.Range("E12").Select()

[void]$wsChangeDemand.Activate()
[void]$wsChangeDemand.Range("E12").Select()
